$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3577.7
$ws.Range("I106").Value = 3307.889
$ws.Range("K106").Value = 3307.889
$ws.Range("M106").Value = -2676.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H122").Value = 4895
$ws.Range("I122").Value = 4895
$ws.Range("K122").Value = 14685
$ws.Range("M122").Value = -12235

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1180.7241
$ws.Range("I5").Value = 1021.3333
$ws.Range("K5").Value = 3063.9999
$ws.Range("M5").Value = -2951.9999
$ws.Range("H135").Value = 1180.7241
$ws.Range("I135").Value = 1021.3333
$ws.Range("K135").Value = 9191.9997
$ws.Range("M135").Value = -6656.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("H126").ClearContents()
$ws.Range("I126").ClearContents()
$ws.Range("J126").ClearContents()
$ws.Range("K126").ClearContents()
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("H136").ClearContents()
$ws.Range("I136").ClearContents()
$ws.Range("J136").ClearContents()
$ws.Range("K136").ClearContents()
$ws.Range("L136").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("M140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4293.3335
$ws.Range("J46").Value = 1999
$ws.Range("L46").Value = 1999
$ws.Range("N46").Value = -2375
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -41308
$ws.Range("H124").ClearContents()
$ws.Range("I124").ClearContents()
$ws.Range("J124").ClearContents()
$ws.Range("K124").ClearContents()
$ws.Range("L124").ClearContents()
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("M131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H136").ClearContents()
$ws.Range("I136").ClearContents()
$ws.Range("J136").ClearContents()
$ws.Range("K136").ClearContents()
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 55000
$ws.Range("J114").Value = 55000
$ws.Range("L114").Value = 55000
$ws.Range("N114").Value = -63678
$ws.Range("H126").Value = 4148.6294
$ws.Range("I126").Value = 2863.8948
$ws.Range("J126").Value = 7199.875
$ws.Range("K126").Value = 8591.6844
$ws.Range("L126").Value = 21599.625
$ws.Range("M126").Value = -26539.625
